# Apply cell updates per the commit diff.
# Cells whose new value looks like a plain number (e.g. "8.00", "0.999")
# must be forced to remain text, matching the original inlineStr/text cells,
# otherwise Excel auto-converts them to numeric values (losing formatting like
# trailing zeros, and breaking the "t=inlineStr"/string nature of the cell).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.550.17"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "3.443.74"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.97"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.92"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.484"
$ws.Range("E8").Value = "  +0.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "8.04"
$ws.Range("E9").Value = "  +5.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.123"
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("E11").Value = "  +2.46%  "
$ws.Range("D12").Value = "4.032.85"
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("E13").Value = "  +1.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.28"
$ws.Range("E14").Value = "  -5.25%  "
$ws.Range("D15").Value = "3.449.26"
$ws.Range("E15").Value = "  -0.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000171"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").Value = "62.589.86"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.34"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.65"
$ws.Range("E19").Value = "  +1.45%  "
$ws.Range("E20").Value = "  -3.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "386.08"
$ws.Range("E21").Value = "  -0.90%  "
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "75.11"
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").Value = "3.582.46"
$ws.Range("E25").Value = "  -1.06%  "
$ws.Range("E26").Value = "  -1.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.182"
$ws.Range("E27").Value = "  +0.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.61"
$ws.Range("E28").Value = "  -0.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.00"
$ws.Range("E30").Value = "  -2.57%  "
$ws.Range("E31").Value = "  -1.75%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("E33").Value = "  -4.78%  "
$ws.Range("E34").Value = "  -2.54%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.37"
$ws.Range("E35").Value = "  +1.61%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.63"
$ws.Range("E36").Value = "  +4.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "31.83"
$ws.Range("E37").Value = "  +1.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.93"
$ws.Range("E38").Value = "  -2.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "169.02"
$ws.Range("E39").Value = "  -1.21%  "
$ws.Range("D40").Value = "3.477.65"
$ws.Range("E40").Value = "  -0.98%  "
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.784"
$ws.Range("E42").Value = "  -2.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.58"
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.37"
$ws.Range("E44").Value = "  -3.00%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.68"
$ws.Range("E45").Value = "  -1.55%  "
$ws.Range("E46").Value = "  -1.93%  "
$ws.Range("D47").Value = "2.573.48"
$ws.Range("E47").Value = "  -1.04%  "
$ws.Range("E48").Value = "  +1.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.23"
$ws.Range("E49").Value = "  -1.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.59"
$ws.Range("E50").Value = "  -4.34%  "
$ws.Range("E51").Value = "  -0.06%  "
